$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.746.20'
$ws.Range('E2').Value = '  +2.26%  '
$ws.Range('D3').Value = '1.692.06'
$ws.Range('E3').Value = '  +3.25%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '221.95'
$ws.Range('E5').Value = '  +2.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.522'
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '30.80'
$ws.Range('E8').Value = '  +4.31%  '
$ws.Range('E9').Value = '  +1.72%  '
$ws.Range('E10').Value = '  +2.04%  '
$ws.Range('E11').Value = '  -1.49%  '
$ws.Range('D12').Value = '1.937.59'
$ws.Range('E12').Value = '  +3.52%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.63'
$ws.Range('E13').Value = '  +10.34%  '
$ws.Range('E14').Value = '  +7.80%  '
$ws.Range('D15').Value = '1.700.13'
$ws.Range('E15').Value = '  +3.36%  '
$ws.Range('E16').Value = '  +2.50%  '
$ws.Range('D17').Value = '30.781.26'
$ws.Range('E17').Value = '  +2.41%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '66.48'
$ws.Range('E18').Value = '  +2.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '249.13'
$ws.Range('E19').Value = '  -0.15%  '
$ws.Range('E20').Value = '  +1.56%  '
$ws.Range('E21').Value = '  -0.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.27'
$ws.Range('E22').Value = '  +3.57%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.30'
$ws.Range('E23').Value = '  +2.03%  '
$ws.Range('E24').Value = '  +2.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.75'
$ws.Range('E25').Value = '  -1.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.94'
$ws.Range('E26').Value = '  +1.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.111'
$ws.Range('E27').Value = '  -0.19%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.75'
$ws.Range('E28').Value = '  +0.92%  '
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0501'
$ws.Range('E30').Value = '  +1.58%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.15'
$ws.Range('E31').Value = '  +0.96%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.48'
$ws.Range('E32').Value = '  +1.36%  '
$ws.Range('D33').Value = '1.518.00'
$ws.Range('E33').Value = '  +5.44%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.31'
$ws.Range('E34').Value = '  +2.38%  '
$ws.Range('E35').Value = '  +4.27%  '
$ws.Range('E36').Value = '  -0.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '81.76'
$ws.Range('E37').Value = '  +7.24%  '
$ws.Range('E38').Value = '  +4.66%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.586'
$ws.Range('E39').Value = '  +4.44%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.70'
$ws.Range('E40').Value = '  -5.92%  '
$ws.Range('E41').Value = '  +1.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.854'
$ws.Range('E42').Value = '  +1.52%  '
$ws.Range('E43').Value = '  +0.94%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0505'
$ws.Range('E44').Value = '  +1.19%  '
$ws.Range('E45').Value = '  -1.39%  '
$ws.Range('E46').Value = '  -0.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '52.42'
$ws.Range('E47').Value = '  -4.90%  '
$ws.Range('E49').Value = '  +0.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '95.65'
$ws.Range('E50').Value = '  +5.53%  '
$ws.Range('D51').Value = '0.0₆0116'
$ws.Range('E51').Value = '  +4.92%  '
